# Applies the cryptos-list price/volume refresh described in the commit.
# D-column prices are forced to remain text (leading apostrophe, as Excel
# does for literal/quoted text) so values like "0.999" or "6.70" are not
# auto-coerced into numbers and keep their original formatting exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.228.69"
$ws.Range("E2").Value = "  -4.76%  "
$ws.Range("D3").Value = "'3.310.79"
$ws.Range("E3").Value = "  -5.32%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'564.49"
$ws.Range("E5").Value = "  -4.04%  "
$ws.Range("D6").Value = "'127.31"
$ws.Range("E6").Value = "  -4.69%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'3.309.15"
$ws.Range("E8").Value = "  -5.38%  "
$ws.Range("D9").Value = "'0.479"
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").Value = "'7.34"
$ws.Range("E10").Value = "  -4.34%  "
$ws.Range("E11").Value = "  -4.66%  "
$ws.Range("D12").Value = "'0.375"
$ws.Range("E12").Value = "  -2.77%  "
$ws.Range("D13").Value = "'3.868.69"
$ws.Range("E13").Value = "  -5.45%  "
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "'3.305.08"
$ws.Range("E15").Value = "  -5.47%  "
$ws.Range("D16").Value = "'0.0000168"
$ws.Range("E16").Value = "  -6.33%  "
$ws.Range("D17").Value = "'24.71"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "'61.242.73"
$ws.Range("E18").Value = "  -4.66%  "
$ws.Range("D19").Value = "'13.51"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "'5.65"
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("D21").Value = "'9.04"
$ws.Range("E21").Value = "  -9.61%  "
$ws.Range("D22").Value = "'353.51"
$ws.Range("E22").Value = "  -8.41%  "
$ws.Range("D23").Value = "'0.555"
$ws.Range("E23").Value = "  -3.92%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'3.436.92"
$ws.Range("E25").Value = "  -5.46%  "
$ws.Range("D26").Value = "'69.24"
$ws.Range("E26").Value = "  -6.81%  "
$ws.Range("D27").Value = "'0.0000107"
$ws.Range("E27").Value = "  -7.18%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").Value = "'7.15"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("D30").Value = "'7.85"
$ws.Range("E30").Value = "  -3.31%  "
$ws.Range("D31").Value = "'1.41"
$ws.Range("E31").Value = "  -5.00%  "
$ws.Range("D32").Value = "'2.10"
$ws.Range("E32").Value = "  -6.42%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").Value = "'0.149"
$ws.Range("E34").Value = "  -3.45%  "
$ws.Range("D35").Value = "'3.334.18"
$ws.Range("E35").Value = "  -5.40%  "
$ws.Range("D36").Value = "'22.53"
$ws.Range("E36").Value = "  -3.09%  "
$ws.Range("D37").Value = "'5.23"
$ws.Range("E37").Value = "  -2.82%  "
$ws.Range("D38").Value = "'6.79"
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("D39").Value = "'160.11"
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("E40").Value = "  -4.11%  "
$ws.Range("D41").Value = "'0.0758"
$ws.Range("E41").Value = "  -3.18%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "'41.06"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("D44").Value = "'4.36"
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("D45").Value = "'0.743"
$ws.Range("E45").Value = "  -7.84%  "
$ws.Range("E46").Value = "  -5.53%  "
$ws.Range("E47").Value = "  -5.53%  "
$ws.Range("D48").Value = "'22.28"
$ws.Range("E48").Value = "  -8.23%  "
$ws.Range("D49").Value = "'6.70"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").Value = "'0.863"
$ws.Range("E50").Value = "  -5.95%  "
$ws.Range("D51").Value = "'21.02"
$ws.Range("E51").Value = "  +0.70%  "
